$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-03"

# Update the "May (through 05-02)" label to "May (through 05-03)"
$ws.Range("A6").Value = "May (through 05-03)"

# Update the May row (row 6) with the new data point added on 2022-05-11
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 5
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 10

# Update the Total row (row 7) to reflect the new May figures
$ws.Range("C7").Value = 163
$ws.Range("D7").Value = 258
$ws.Range("G7").Value = 265
$ws.Range("H7").Value = 532
$ws.Range("I7").Value = 561
